$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 340, shifting the existing
# rows 340:363 down to 341:364 (same as Excel's right-click > Insert).
$ws.Rows.Item(340).Insert()

# Populate the newly inserted row 340 with the new weekly record.
$ws.Cells.Item(340, 1).Value = 9
$ws.Cells.Item(340, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(340, 3).Value = "Metropolitana"
$ws.Cells.Item(340, 4).Value = 44931
$ws.Cells.Item(340, 5).Value = 13
$ws.Cells.Item(340, 6).Value = 300000001
$ws.Cells.Item(340, 7).Value = "Rabanito"
$ws.Cells.Item(340, 8).Value = "Sin especificar"
$ws.Cells.Item(340, 9).Value = "Primera"
$ws.Cells.Item(340, 10).Value = 7000
$ws.Cells.Item(340, 11).Value = 3000
$ws.Cells.Item(340, 12).Value = 3000
$ws.Cells.Item(340, 13).Value = 3000
$ws.Cells.Item(340, 14).Value = "$/cien unidades (volumen en unidades)"
$ws.Cells.Item(340, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(340, 16).Value = 30
$ws.Cells.Item(340, 17).Value = 100
$ws.Cells.Item(340, 18).Value = "Hortaliza"

# Make sure the new D340 keeps the same date-style formatting (s="2")
# as the rest of the Fecha column.
$ws.Cells.Item(340, 4).NumberFormat = $ws.Cells.Item(341, 4).NumberFormat
